# Regenerate the handoff/handback status report.
#
# The underlying data for the two tracked source files swapped places:
#   - bed71e79-5634-4f78-be6c-f2c1dca6cf5b  is now the "top" row (still
#     handed back / in sync with en-US, but with a refreshed timestamp).
#   - 0055612a-0e2e-4e43-b43c-dc47a2018e97  moved to the second row and is
#     now "Ready for handoff" (it had round-tripped back out of sync).
#
# Hyperlink targets (the rIds / actual URLs) do not change - only the
# *values* shown in the grid and the hyperlinks' visible display text do,
# so we update cell values and then walk each sheet's Hyperlinks collection
# (via foreach, which yields live objects whose TextToDisplay can be set
# in place) to resync the display text without touching styles or rIds.

$wb = $excel.ActiveWorkbook

function Set-LinkText($ws, $addr, $text) {
    foreach ($h in $ws.Hyperlinks()) {
        if ($h.Range().Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-03-30 11:02:32"

$ov.Range("A3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-30 11:02:32"

Set-LinkText $ov '$A$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
Set-LinkText $ov '$A$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-30 11:02:18"
$zh.Range("F2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$zh.Range("G2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"

$zh.Range("A3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-30 11:02:18"
$zh.Range("F3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$zh.Range("G3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"

Set-LinkText $zh '$A$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
Set-LinkText $zh '$D$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
Set-LinkText $zh '$F$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
Set-LinkText $zh '$G$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
Set-LinkText $zh '$A$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
Set-LinkText $zh '$D$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"
Set-LinkText $zh '$F$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
Set-LinkText $zh '$G$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
$de.Range("E2").Value = "2016-03-30 11:02:32"
$de.Range("F2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$de.Range("G2").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"

$de.Range("A3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
$de.Range("E3").Value = "2016-03-30 11:02:32"
$de.Range("F3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$de.Range("G3").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"

Set-LinkText $de '$A$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
Set-LinkText $de '$D$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
Set-LinkText $de '$F$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
Set-LinkText $de '$G$2' "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
Set-LinkText $de '$A$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
Set-LinkText $de '$D$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
Set-LinkText $de '$F$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
Set-LinkText $de '$G$3' "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
